$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = -7.546000000000001
$ws.Range("A3").Value = -21.784
$ws.Range("E3").Value = 16.572
$ws.Range("E12").Value = 17.429
$ws.Range("A14").Value = -21.659
$ws.Range("A16").Value = -22.075
$ws.Range("D18").Value = -8.258999999999999
$ws.Range("A21").Value = -20.034
$ws.Range("A23").Value = -20.302
$ws.Range("D24").Value = -7.591999999999999
$ws.Range("E24").Value = 16.873
$ws.Range("A25").Value = -21.608
$ws.Range("D25").Value = -7.281000000000001
$ws.Range("E25").Value = 16.935
$ws.Range("A26").Value = -21.122
$ws.Range("D27").Value = -7.922000000000001
$ws.Range("A29").Value = -21.232
$ws.Range("D30").Value = -7.305
$ws.Range("D31").Value = -7.888999999999998
$ws.Range("D39").Value = -7.502
$ws.Range("A40").Value = -20.197
$ws.Range("E41").Value = 16.484
$ws.Range("D48").Value = -7.101000000000001
$ws.Range("D51").Value = -8.326000000000001
$ws.Range("D52").Value = -8.083
$ws.Range("A53").Value = -21.938
$ws.Range("E53").Value = 16.484
$ws.Range("D55").Value = -8.179
$ws.Range("D56").Value = -8.331999999999999
$ws.Range("E56").Value = 16.36
$ws.Range("A57").Value = -22.606
$ws.Range("D57").Value = -8.217000000000002
$ws.Range("E57").Value = 16.481
$ws.Range("E58").Value = 16.478
$ws.Range("A59").Value = -22.571
$ws.Range("D60").Value = -8.408000000000001
$ws.Range("E61").Value = 16.652
$ws.Range("E63").Value = 17.392
$ws.Range("E64").Value = 17.186
$ws.Range("A65").Value = -21.482
$ws.Range("A69").Value = -21.52
$ws.Range("E70").Value = 17.464
$ws.Range("E72").Value = 16.97000000000001
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("D74").Value = -7.890000000000001
$ws.Range("A79").Value = -20.896
$ws.Range("A83").Value = -21.919
$ws.Range("E86").Value = 16.362
$ws.Range("D89").Value = -6.753
$ws.Range("E89").Value = 17.42
$ws.Range("D90").Value = -7.475999999999999
$ws.Range("A91").Value = -21.527
$ws.Range("D92").Value = -6.701000000000001
$ws.Range("A93").Value = -21.216
$ws.Range("E98").Value = 16.491
$ws.Range("A100").Value = -22.111
$ws.Range("E100").Value = 16.562
$ws.Range("E102").Value = 16.49
